# Apply "correções conf. feedback Professor" edit:
# Append a trailing period to the risk "Nome" (column C) for rows 3, 4, 5, 6 and 8,
# then leave the selection on C8 with the sheet scrolled so row 4 is the first visible row
# (mirrors the final editing position left by the author).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Riscos")

$ws.Range("C3").Value = "Problema na contratação do Data Center."
$ws.Range("C4").Value = "Programador com domínio pleno nas duas linguagens principais PHP e Java utilizadas no sistema."
$ws.Range("C6").Value = "Mudança na expectativa dos Stakeholder´s quanto às funcionalidades do Sistema de Rastreamento."
$ws.Range("C5").Value = "Realizar correções na aplicação principal não previstas e muito frequentes."
$ws.Range("C8").Value = "Queima de equipamentos diversos que serão utilizados pela a equipe de desenvolvimento de software."

$ws.Activate()
$ws.Range("C8").Select()
$excel.ActiveWindow.ScrollRow = 4
